$wb = $excel.ActiveWorkbook

# --- Benchmarks sheet: selection moves from D4 to C23 (sheet stays inactive) ---
$wsBenchmarks = $wb.Worksheets.Item("Benchmarks")
$wsBenchmarks.Range("C23").Select()

# --- Config1 sheet: fill in LRU (row43) + finish hmmer rows 41/42 for lbm benchmark ---
$wsConfig1 = $wb.Worksheets.Item("Config1")
$wsConfig1.Activate()

# Row 41 (hmmer / Hawkeye)
$wsConfig1.Range("C41").Value = 50000000
$wsConfig1.Range("D41").Value = 44574971
$wsConfig1.Range("E41").Value = 328304
$wsConfig1.Range("F41").Value = 318182
$wsConfig1.Range("G41").Value = 10122

# Row 42 (hmmer / OPTGen)
$wsConfig1.Range("C42").Value = 50000000
$wsConfig1.Range("D42").Value = 44574971
$wsConfig1.Range("E42").Value = 5252
$wsConfig1.Range("F42").Value = 3689
$wsConfig1.Range("G42").Formula = "=E42-F42"
$wsConfig1.Range("J42").Formula = "=F42/E42"

# Row 43 (lbm / LRU)
$wsConfig1.Range("C43").Value = 50000000
$wsConfig1.Range("D43").Value = 93261366
$wsConfig1.Range("E43").Value = 2641894
$wsConfig1.Range("F43").Value = 1161101
$wsConfig1.Range("G43").Value = 1480793

$wsConfig1.Range("C44").Select()

# --- Config2 sheet: same rows, different simulation numbers ---
$wsConfig2 = $wb.Worksheets.Item("Config2")
$wsConfig2.Activate()

# Row 41 (hmmer / Hawkeye)
$wsConfig2.Range("C41").Value = 50000000
$wsConfig2.Range("D41").Value = 43537224
$wsConfig2.Range("E41").Value = 33170
$wsConfig2.Range("F41").Value = 270148
$wsConfig2.Range("G41").Value = 61222

# Row 42 (hmmer / OPTGen)
$wsConfig2.Range("C42").Value = 50000000
$wsConfig2.Range("D42").Value = 43537224
$wsConfig2.Range("E42").Value = 1797
$wsConfig2.Range("F42").Value = 1620
$wsConfig2.Range("G42").Value = 177
$wsConfig2.Range("J42").Formula = "=F42/E42"

# Row 43 (lbm / LRU)
$wsConfig2.Range("C43").Value = 50000000
$wsConfig2.Range("D43").Value = 80458099
$wsConfig2.Range("E43").Value = 2664603
$wsConfig2.Range("F43").Value = 1183779
$wsConfig2.Range("G43").Value = 1480824

$wsConfig2.Range("C44").Select()

# Re-activate Config1 so it matches the workbook's active tab / selection
$wsConfig1.Activate()
$wsConfig1.Range("C44").Select()

Write-Output "done"
